# Updates the cryptos list: new Price (column D) and Volume(1h) (column E)
# values for each coin row, matching the latest scrape from GitHub Actions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells whose new text would otherwise be auto-parsed by Excel as a number
# (single decimal point, all digits) get written with a leading apostrophe
# to force text, then have their style reset to Normal so no stray
# quote-prefix/number-format style is left behind on the cell.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '30.085.25'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue 'D5' '346.69'
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('E6').Value = '  +0.00%  '
Set-TextValue 'D7' '0.5198'
$ws.Range('E7').Value = '  +0.25%  '
Set-TextValue 'D8' '0.4480'
$ws.Range('E8').Value = '  +0.72%  '
Set-TextValue 'D9' '54.05'
$ws.Range('E9').Value = '  +3.09%  '
Set-TextValue 'D10' '0.09398'
$ws.Range('E10').Value = '  -0.93%  '
$ws.Range('E11').Value = '  +0.42%  '
Set-TextValue 'D12' '25.46'
$ws.Range('E12').Value = '  +0.22%  '
Set-TextValue 'D13' '8.693'
$ws.Range('E13').Value = '  +7.25%  '
Set-TextValue 'D14' '6.989'
$ws.Range('E14').Value = '  +3.60%  '
$ws.Range('D15').Value = '2.092.63'
$ws.Range('E15').Value = '  -0.68%  '
Set-TextValue 'D16' '102.69'
$ws.Range('E16').Value = '  +3.08%  '
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('E19').Value = '  +4.33%  '
Set-TextValue 'D20' '0.06705'
$ws.Range('E20').Value = '  +0.09%  '
Set-TextValue 'D21' '6.314'
$ws.Range('E21').Value = '  +1.85%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').Value = '30.113.71'
$ws.Range('E23').Value = '  -0.15%  '
Set-TextValue 'D24' '12.74'
$ws.Range('E24').Value = '  -0.03%  '
Set-TextValue 'D25' '2.334'
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('D26').Value = '2.361.64'
$ws.Range('E26').Value = '  +0.33%  '
Set-TextValue 'D27' '22.22'
$ws.Range('E27').Value = '  +0.67%  '
Set-TextValue 'D28' '2.555'
$ws.Range('E28').Value = '  +0.63%  '
Set-TextValue 'D29' '162.81'
Set-TextValue 'D30' '134.23'
$ws.Range('E30').Value = '  +0.39%  '
Set-TextValue 'D31' '1.161'
$ws.Range('E31').Value = '  +0.24%  '
Set-TextValue 'D32' '1.790'
$ws.Range('E32').Value = '  +9.71%  '
Set-TextValue 'D33' '0.1058'
$ws.Range('E33').Value = '  +0.10%  '
Set-TextValue 'D34' '6.293'
$ws.Range('E34').Value = '  +0.57%  '
Set-TextValue 'D35' '6.683'
$ws.Range('E35').Value = '  +8.24%  '
Set-TextValue 'D36' '3.971'
$ws.Range('E36').Value = '  +0.69%  '
Set-TextValue 'D37' '10.80'
$ws.Range('E37').Value = '  +6.31%  '
Set-TextValue 'D38' '0.02644'
$ws.Range('E38').Value = '  +2.51%  '
Set-TextValue 'D39' '0.06890'
$ws.Range('E39').Value = '  +1.50%  '
Set-TextValue 'D40' '0.7142'
$ws.Range('E40').Value = '  +2.57%  '
Set-TextValue 'D41' '12.75'
$ws.Range('E41').Value = '  +1.79%  '
Set-TextValue 'D42' '0.2257'
$ws.Range('E42').Value = '  -1.30%  '
Set-TextValue 'D43' '1.330'
$ws.Range('E43').Value = '  +1.98%  '
Set-TextValue 'D44' '0.6932'
$ws.Range('E44').Value = '  +3.46%  '
Set-TextValue 'D45' '14.71'
$ws.Range('E45').Value = '  +3.00%  '
Set-TextValue 'D46' '2.400'
$ws.Range('E46').Value = '  +4.73%  '
$ws.Range('E47').Value = '  +0.03%  '
Set-TextValue 'D48' '3.635'
$ws.Range('E48').Value = '  -0.04%  '
Set-TextValue 'D49' '1.262'
$ws.Range('E49').Value = '  +7.43%  '
Set-TextValue 'D50' '0.00000000357'
$ws.Range('E50').Value = '  -0.47%  '
